# Atualização Sprint backlog da semana
#
# - Adds a new "Sprint Backlog (25-04)" worksheet (a continuation of the
#   "Sprint Backlog (18-04)" backlog table, now with 10 tasks instead of 7).
# - Shrinks a few row heights on "Sprint Backlog (18-04)" and clears its
#   stale selection / tab-active state now that the new sheet takes over.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122
$xlLeft = -4131
$xlCenter = -4108
$xlContinuous = 1
$xlMedium = -4138
$xlNone = -4142
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

# ---------------------------------------------------------------------
# 1. Duplicate "Sprint Backlog (18-04)" -> "Sprint Backlog (25-04)"
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("Sprint Backlog (18-04)")
$srcSheet.Copy([System.Reflection.Missing]::Value, $srcSheet)
$ws = $wb.Worksheets.Item("Sprint Backlog (18-04) (2)")
$ws.Name = "Sprint Backlog (25-04)"

# Drop the 3 small merges inherited from the copy - the new layout uses a
# single B2:B9 merge instead.
$ws.Range("B2:B4").UnMerge()
$ws.Range("B5:B6").UnMerge()
$ws.Range("B7:B8").UnMerge()

# ---------------------------------------------------------------------
# 2. Column widths for the new sheet
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 19
$ws.Columns.Item(3).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 15.7109375

Write-Host "stage1 ok"
